# Updates the cryptocurrency price/volume table to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.430.03"
$ws.Range("E2").Value = "  +4.67%  "
$ws.Range("D3").Value = "2.488.07"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.40"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.20"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.524"
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.76"
$ws.Range("E10").Value = "  +5.91%  "
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.13"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "2.875.40"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "2.500.31"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "47.303.71"
$ws.Range("E18").Value = "  +4.63%  "
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.55"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.59"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.39"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.37"
$ws.Range("E24").Value = "  +5.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +2.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.13"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.99"
$ws.Range("E28").Value = "  +4.27%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.95"
$ws.Range("E30").Value = "  +6.04%  "
$ws.Range("E31").Value = "  +6.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.46"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.82"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0779"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.94"
$ws.Range("E37").Value = "  +3.26%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  +3.54%  "
$ws.Range("E39").Value = "  +4.39%  "
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "121.48"
$ws.Range("E42").Value = "  -3.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.12"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0296"
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("D45").Value = "1.960.68"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.19"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.80"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("E50").Value = "  +11.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.29"
$ws.Range("E51").Value = "  +3.71%  "

Write-Output "Applied cryptos update"
